# Rename the two worksheets to match the excel-importer's expected sheet
# names (Sheet1 -> Sheet4, Sheet2 -> Sheet5). Sheet1/Sheet2 correspond to
# rId1/rId2 i.e. Worksheets.Item(1)/Item(2) respectively.
$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item(1)
$ws4.Name = "Sheet4"

$ws5 = $wb.Worksheets.Item(2)
$ws5.Name = "Sheet5"

# --- Sheet4 (was Sheet1): scroll the view further down, from A1974 to
# A1986. No cell should end up selected other than the default A1 (the
# sheet keeps no <selection> override), just the visible window moves.
$ws4.Activate()
$winTop = $excel.ActiveWindow
$winTop.ScrollRow = 1986
$winTop.ScrollColumn = 1

# --- Sheet5 (was Sheet2): this is the active/visible tab. Its viewport
# scrolls down so row 1974 is the top-left visible row, and the selected
# cell moves from E11 down to L1984 (which is what produces/justifies the
# scroll in the first place).
$ws5.Activate()
[void]$ws5.Range("L1984").Select()
$winBottom = $excel.ActiveWindow
$winBottom.ScrollRow = 1974
$winBottom.ScrollColumn = 1

# Leave Sheet5 as the active/selected tab (matches activeTab/tabSelected
# staying on the second sheet before and after the edit).
$ws5.Activate()
